$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 70
$ws.Range("H70").Value = 1503.1818
$ws.Range("I70").Value = 1426.6
$ws.Range("J70").Value = 1567
$ws.Range("K70").Value = 4279.799999999999
$ws.Range("L70").Value = 4701
$ws.Range("M70").Value = -4009.799999999999
$ws.Range("N70").Value = -5241

# Row 73
$ws.Range("H73").Value = 1503.1818
$ws.Range("I73").Value = 1426.6
$ws.Range("J73").Value = 1567
$ws.Range("K73").Value = 4279.799999999999
$ws.Range("L73").Value = 4701
$ws.Range("M73").Value = -3343.799999999999
$ws.Range("N73").Value = -6573

# Row 138
$ws.Range("H138").Value = 2335.1428
$ws.Range("I138").Value = 3055.0715
$ws.Range("J138").Value = 2047.1714
$ws.Range("K138").Value = 9165.2145
$ws.Range("L138").Value = 6141.5142
$ws.Range("M138").Value = -4025.2145
$ws.Range("N138").Value = -16421.5142

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 594735.2
$ws.Range("I32").Value = 773516.5600000001
$ws.Range("J32").Value = 17134
$ws.Range("K32").Value = 773516.5600000001
$ws.Range("L32").Value = 17134
$ws.Range("M32").Value = -773229.5600000001
$ws.Range("N32").Value = -17708

# Row 45
$ws.Range("H45").Value = 2487.5652
$ws.Range("I45").Value = 1438.4615
$ws.Range("K45").Value = 1438.4615
$ws.Range("M45").Value = -1061.4615

# Row 110
$ws.Range("H110").Value = 2164.2
$ws.Range("I110").Value = 1607
$ws.Range("J110").Value = 3000
$ws.Range("K110").Value = 1607
$ws.Range("L110").Value = 3000
$ws.Range("M110").Value = 438
$ws.Range("N110").Value = -7090

# Row 132
$ws.Range("H132").Value = 3676.0789
$ws.Range("I132").Value = 3390.5908
$ws.Range("J132").Value = 4068.625
$ws.Range("K132").Value = 10171.7724
$ws.Range("L132").Value = 12205.875
$ws.Range("M132").Value = -7641.7724
$ws.Range("N132").Value = -17265.875

$ws = $wb.Worksheets.Item("BSM")
# Row 25
$ws.Range("H25").Value = 32403.2
$ws.Range("I25").Value = 17000
$ws.Range("J25").Value = 42672
$ws.Range("K25").Value = 17000
$ws.Range("L25").Value = 42672
$ws.Range("M25").Value = -16765
$ws.Range("N25").Value = -43142

# Row 105
$ws.Range("H105").Value = 3425
$ws.Range("I105").Value = 3450
$ws.Range("K105").Value = 3450
$ws.Range("M105").Value = -1703

# Row 107
$ws.Range("H107").Value = 2972.8333
$ws.Range("I107").Value = 1905.5
$ws.Range("J107").Value = 3506.5
$ws.Range("K107").Value = 1905.5
$ws.Range("L107").Value = 3506.5
$ws.Range("M107").Value = 14.5
$ws.Range("N107").Value = -7346.5

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1600
$ws.Range("I16").Value = 1000
$ws.Range("J16").Value = 2200
$ws.Range("K16").Value = 1000
$ws.Range("L16").Value = 2200
$ws.Range("M16").Value = -713
$ws.Range("N16").Value = -2774

# Row 113
$ws.Range("H113").Value = 1600
$ws.Range("I113").Value = 1000
$ws.Range("J113").Value = 2200
$ws.Range("K113").Value = 1000
$ws.Range("L113").Value = 2200
$ws.Range("M113").Value = 1170
$ws.Range("N113").Value = -6540

$ws = $wb.Worksheets.Item("CUL")
# Row 17
$ws.Range("H17").Value = 578.5714
$ws.Range("I17").Value = 600
$ws.Range("K17").Value = 1800
$ws.Range("M17").Value = -1631

# Row 40
$ws.Range("H40").Value = 219.47058
$ws.Range("I40").Value = 181.7
$ws.Range("J40").Value = 273.42856
$ws.Range("K40").Value = 726.8
$ws.Range("L40").Value = 1093.71424
$ws.Range("M40").Value = -657.8
$ws.Range("N40").Value = -1231.71424

# Row 122
$ws.Range("H122").Value = 6267.278
$ws.Range("I122").Value = 524.2308
$ws.Range("K122").Value = 4718.077200000001
$ws.Range("M122").Value = -2268.077200000001

# Row 138
$ws.Range("H138").Value = 5035.2856
$ws.Range("I138").Value = 1346
$ws.Range("J138").Value = 7084.8887
$ws.Range("K138").Value = 4038
$ws.Range("L138").Value = 21254.6661
$ws.Range("M138").Value = 1102
$ws.Range("N138").Value = -31534.6661

# Row 140
$ws.Range("H140").Value = 1705.8846
$ws.Range("I140").Value = 981.7368
$ws.Range("J140").Value = 3671.4285
$ws.Range("K140").Value = 2945.2104
$ws.Range("L140").Value = 11014.2855
$ws.Range("M140").Value = 2234.7896
$ws.Range("N140").Value = -21374.2855

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 6149.9546
$ws.Range("I70").Value = 6571.2856
$ws.Range("J70").Value = 5953.3335
$ws.Range("K70").Value = 6571.2856
$ws.Range("L70").Value = 5953.3335
$ws.Range("M70").Value = -6301.2856
$ws.Range("N70").Value = -6493.3335

# Row 73
$ws.Range("H73").Value = 6149.9546
$ws.Range("I73").Value = 6571.2856
$ws.Range("J73").Value = 5953.3335
$ws.Range("K73").Value = 6571.2856
$ws.Range("L73").Value = 5953.3335
$ws.Range("M73").Value = -5635.2856
$ws.Range("N73").Value = -7825.3335

# Row 113
$ws.Range("H113").Value = 2942
$ws.Range("I113").Value = 2913
$ws.Range("J113").Value = 2956.5
$ws.Range("K113").Value = 2913
$ws.Range("L113").Value = 2956.5
$ws.Range("M113").Value = -743
$ws.Range("N113").Value = -7296.5

# Row 132
$ws.Range("H132").Value = 2886.5
$ws.Range("I132").Value = 2160.5557
$ws.Range("K132").Value = 6481.6671
$ws.Range("M132").Value = -3951.6671

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 10305
$ws.Range("I22").Value = 1467
$ws.Range("J22").Value = 13251
$ws.Range("K22").Value = 1467
$ws.Range("L22").Value = 13251
$ws.Range("M22").Value = -1172
$ws.Range("N22").Value = -13841

# Row 27
$ws.Range("H27").Value = 10305
$ws.Range("I27").Value = 1467
$ws.Range("J27").Value = 13251
$ws.Range("K27").Value = 1467
$ws.Range("L27").Value = 13251
$ws.Range("M27").Value = -1360
$ws.Range("N27").Value = -13465

# Row 61
$ws.Range("H61").Value = 4236.6
$ws.Range("I61").Value = 4555.4443
$ws.Range("J61").Value = 3758.3333
$ws.Range("K61").Value = 4555.4443
$ws.Range("L61").Value = 3758.3333
$ws.Range("M61").Value = -4353.4443
$ws.Range("N61").Value = -4162.3333

# Row 110
$ws.Range("H110").Value = 50000
$ws.Range("J110").Value = 50000
$ws.Range("L110").Value = 50000
$ws.Range("M110").Value = -58180

# Row 113
$ws.Range("H113").Value = 4236.6
$ws.Range("I113").Value = 4555.4443
$ws.Range("J113").Value = 3758.3333
$ws.Range("K113").Value = 4555.4443
$ws.Range("L113").Value = 3758.3333
$ws.Range("M113").Value = -2385.4443
$ws.Range("N113").Value = -8098.3333

# Row 119
$ws.Range("H119").Value = 48000
$ws.Range("J119").Value = 48000
$ws.Range("L119").Value = 48000
$ws.Range("N119").Value = -57676

$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 682.5714
$ws.Range("I107").Value = 719.5
$ws.Range("J107").Value = 633.3333
$ws.Range("K107").Value = 2158.5
$ws.Range("L107").Value = 1899.9999
$ws.Range("M107").Value = -238.5
$ws.Range("N107").Value = -5739.9999
